$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-04 Monday" "2024-03-05 Tuesday"

Replace-Text "300÷3=" "474÷4="
Replace-Text "517÷8=" "429÷2="
Replace-Text "133÷3=" "628÷7="
Replace-Text "199÷3=" "580÷5="
Replace-Text "858÷6=" "220÷7="
Replace-Text "770÷6=" "301÷2="
Replace-Text "392÷3=" "648÷7="
Replace-Text "805÷7=" "947÷8="
Replace-Text "141÷8=" "360÷9="
Replace-Text "682÷2=" "411÷9="
Replace-Text "820÷4=" "717÷9="
Replace-Text "984÷6=" "416÷2="
Replace-Text "410÷5=" "918÷8="
Replace-Text "900÷2=" "609÷8="
Replace-Text "608÷5=" "452÷3="
Replace-Text "314÷6=" "860÷2="
Replace-Text "524÷2=" "665÷6="
Replace-Text "940÷4=" "786÷8="
Replace-Text "350÷2=" "335÷9="
Replace-Text "329÷5=" "447÷8="
Replace-Text "954÷6=" "539÷4="
Replace-Text "256÷8=" "189÷6="
Replace-Text "606÷6=" "464÷7="
Replace-Text "224÷6=" "381÷7="
Replace-Text "771÷6=" "495÷7="
